# Update "想去人数" (want-to-go count) figures in both the "展览" sheet
# and the "全部类型" sheet to reflect the refreshed data pull.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 69
$ws1.Range("F4").Value = 1493
$ws1.Range("F5").Value = 577
$ws1.Range("F6").Value = 1062
$ws1.Range("F7").Value = 11029
$ws1.Range("F8").Value = 11029
$ws1.Range("F11").Value = 316
$ws1.Range("F12").Value = 1064
$ws1.Range("F14").Value = 12210
$ws1.Range("F15").Value = 12737
$ws1.Range("F17").Value = 123
$ws1.Range("F22").Value = 34

# --- Sheet "全部类型" (all categories) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 69
$ws4.Range("F5").Value = 1493
$ws4.Range("F6").Value = 577
$ws4.Range("F7").Value = 1062
$ws4.Range("F8").Value = 11029
$ws4.Range("F9").Value = 11029
$ws4.Range("F12").Value = 316
$ws4.Range("F13").Value = 1064
$ws4.Range("F15").Value = 12210
$ws4.Range("F16").Value = 12737
$ws4.Range("F18").Value = 123
$ws4.Range("F23").Value = 34
